$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric stat corrections (data was off by one day; values recomputed) ---
$ws.Range("AU2").Value = 19
$ws.Range("AV2").Value = 17
$ws.Range("AO4").Value = 21
$ws.Range("AW4").Value = 18
$ws.Range("D5").Value = 51
$ws.Range("E5").Value = 20
$ws.Range("G5").Value = 0.392
$ws.Range("L5").Value = 5.5
$ws.Range("N5").Value = 0.347
$ws.Range("Q5").Value = 0.744
$ws.Range("S5").Value = 30
$ws.Range("T5").Value = 43.4
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 21.5
$ws.Range("AB5").Value = 94.59999999999999
$ws.Range("AC5").Value = -3
$ws.Range("AD5").Value = 20
$ws.Range("AE5").Value = 22
$ws.Range("AG5").Value = 22
$ws.Range("AH5").Value = 8
$ws.Range("AS5").Value = 18
$ws.Range("AU5").Value = 18
$ws.Range("AW5").Value = 13
$ws.Range("AZ5").Value = 15
$ws.Range("BA5").Value = 15
$ws.Range("BB5").Value = 24
$ws.Range("AD6").Value = 8
$ws.Range("D7").Value = 52
$ws.Range("F7").Value = 17
$ws.Range("G7").Value = 0.673
$ws.Range("J7").Value = 78.40000000000001
$ws.Range("K7").Value = 0.463
$ws.Range("N7").Value = 0.337
$ws.Range("P7").Value = 25.7
$ws.Range("Q7").Value = 0.82
$ws.Range("R7").Value = 10.7
$ws.Range("S7").Value = 31.8
$ws.Range("U7").Value = 20.2
$ws.Range("AC7").Value = 4.2
$ws.Range("AD7").Value = 8
$ws.Range("AF7").Value = 5
$ws.Range("AG7").Value = 5
$ws.Range("AH7").Value = 26
$ws.Range("AJ7").Value = 27
$ws.Range("AN7").Value = 23
$ws.Range("AR7").Value = 20
$ws.Range("AU7").Value = 21
$ws.Range("AZ7").Value = 21
$ws.Range("AD8").Value = 8
$ws.Range("AN8").Value = 20
$ws.Range("AZ8").Value = 14
$ws.Range("AD9").Value = 8
$ws.Range("AD10").Value = 8
$ws.Range("AS10").Value = 17
$ws.Range("AD11").Value = 8
$ws.Range("AN11").Value = 21
$ws.Range("AV11").Value = 14
$ws.Range("AG12").Value = 21
$ws.Range("AH12").Value = 19
$ws.Range("AD14").Value = 8
$ws.Range("AR14").Value = 19
$ws.Range("AV14").Value = 18
$ws.Range("AD15").Value = 8
$ws.Range("AH15").Value = 11
$ws.Range("AU15").Value = 25
$ws.Range("D16").Value = 51
$ws.Range("F16").Value = 42
$ws.Range("G16").Value = 0.176
$ws.Range("L16").Value = 4.6
$ws.Range("M16").Value = 13.8
$ws.Range("N16").Value = 0.33
$ws.Range("O16").Value = 18.1
$ws.Range("P16").Value = 25.1
$ws.Range("Q16").Value = 0.72
$ws.Range("R16").Value = 9.300000000000001
$ws.Range("S16").Value = 29.1
$ws.Range("T16").Value = 38.4
$ws.Range("V16").Value = 15
$ws.Range("AA16").Value = 21.5
$ws.Range("AB16").Value = 93.09999999999999
$ws.Range("AD16").Value = 20
$ws.Range("AH16").Value = 8
$ws.Range("AN16").Value = 28
$ws.Range("AO16").Value = 18
$ws.Range("AS16").Value = 29
$ws.Range("AU16").Value = 23
$ws.Range("AV16").Value = 16
$ws.Range("BA16").Value = 14
$ws.Range("AN17").Value = 22
$ws.Range("AD18").Value = 20
$ws.Range("AN18").Value = 27
$ws.Range("AS18").Value = 18
$ws.Range("AN19").Value = 26
$ws.Range("AD20").Value = 20
$ws.Range("AE20").Value = 3
$ws.Range("BC20").Value = 5
$ws.Range("AD21").Value = 8
$ws.Range("AH21").Value = 11
$ws.Range("AN21").Value = 24
$ws.Range("BB21").Value = 23
$ws.Range("AJ22").Value = 25
$ws.Range("AU22").Value = 24
$ws.Range("AS23").Value = 28
$ws.Range("AU23").Value = 22
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = 36
$ws.Range("G24").Value = 0.6919999999999999
$ws.Range("J24").Value = 84.3
$ws.Range("K24").Value = 0.492
$ws.Range("L24").Value = 9
$ws.Range("M24").Value = 23.5
$ws.Range("O24").Value = 17.8
$ws.Range("P24").Value = 22.4
$ws.Range("Q24").Value = 0.795
$ws.Range("S24").Value = 32.3
$ws.Range("T24").Value = 40.8
$ws.Range("U24").Value = 27.3
$ws.Range("Z24").Value = 19.5
$ws.Range("AC24").Value = 5.7
$ws.Range("AD24").Value = 8
$ws.Range("AH24").Value = 17
$ws.Range("AO24").Value = 22
$ws.Range("AT24").Value = 22
$ws.Range("AW24").Value = 17
$ws.Range("BA24").Value = 24
$ws.Range("BC24").Value = 6
$ws.Range("AD25").Value = 8
$ws.Range("AD26").Value = 20
$ws.Range("AD27").Value = 20
$ws.Range("AG27").Value = 7
$ws.Range("AJ27").Value = 26
$ws.Range("AD28").Value = 20
$ws.Range("AN28").Value = 25
$ws.Range("AZ28").Value = 13
$ws.Range("AD29").Value = 20
$ws.Range("AH30").Value = 29
$ws.Range("AT30").Value = 23
$ws.Range("AD31").Value = 8
$ws.Range("BA31").Value = 23

# --- Fix Date column formatting: "2-14-2007-08" -> "2008-02-14" (ISO form) ---
# Leading apostrophe forces literal text so Excel does not reinterpret the
# ISO-looking string as a date serial.
$ws.Range("BF2").Value = "'2008-02-14"
$ws.Range("BF3").Value = "'2008-02-14"
$ws.Range("BF4").Value = "'2008-02-14"
$ws.Range("BF5").Value = "'2008-02-14"
$ws.Range("BF6").Value = "'2008-02-14"
$ws.Range("BF7").Value = "'2008-02-14"
$ws.Range("BF8").Value = "'2008-02-14"
$ws.Range("BF9").Value = "'2008-02-14"
$ws.Range("BF10").Value = "'2008-02-14"
$ws.Range("BF11").Value = "'2008-02-14"
$ws.Range("BF12").Value = "'2008-02-14"
$ws.Range("BF13").Value = "'2008-02-14"
$ws.Range("BF14").Value = "'2008-02-14"
$ws.Range("BF15").Value = "'2008-02-14"
$ws.Range("BF16").Value = "'2008-02-14"
$ws.Range("BF17").Value = "'2008-02-14"
$ws.Range("BF18").Value = "'2008-02-14"
$ws.Range("BF19").Value = "'2008-02-14"
$ws.Range("BF20").Value = "'2008-02-14"
$ws.Range("BF21").Value = "'2008-02-14"
$ws.Range("BF22").Value = "'2008-02-14"
$ws.Range("BF23").Value = "'2008-02-14"
$ws.Range("BF24").Value = "'2008-02-14"
$ws.Range("BF25").Value = "'2008-02-14"
$ws.Range("BF26").Value = "'2008-02-14"
$ws.Range("BF27").Value = "'2008-02-14"
$ws.Range("BF28").Value = "'2008-02-14"
$ws.Range("BF29").Value = "'2008-02-14"
$ws.Range("BF30").Value = "'2008-02-14"
$ws.Range("BF31").Value = "'2008-02-14"

# Strip the transient quote-prefix formatting the apostrophe trick added,
# restoring the Date column cells to their original (default) style.
$ws.Range("BF2:BF31").ClearFormats()
